$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old row 5 (its data has been folded into rows 1-4) ---
$ws.Rows.Item(5).Delete()

# --- Row 1: reset to a blank/placeholder "header-ish" record ---
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = "firstname"
$ws.Range("C1").Value = "lastname"
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 0
$ws.Range("F1").Value = 0
$ws.Range("G1").Value = 111111111
$ws.Range("H1").Value = "Blood Type"
$ws.Range("I1").Value = "Gender"

# --- Row 2: fadi badarni ---
$ws.Range("A2").Value = 209315647
$ws.Range("B2").Value = "fadi"
$ws.Range("C2").Value = "badarni"
$ws.Range("D2").Value = 23
$ws.Range("E2").Value = 82
$ws.Range("F2").Value = 184
$ws.Range("G2").Value = 524183083
$ws.Range("H2").Value = "O+"
$ws.Range("I2").Value = "Female"

# --- Row 3: abedalla shiekh ---
$ws.Range("A3").Value = 207527979
$ws.Range("B3").Value = "abedalla"
$ws.Range("C3").Value = "shiekh"
$ws.Range("D3").Value = 22
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 180
$ws.Range("G3").Value = 15241830
$ws.Range("H3").Value = "B"
$ws.Range("I3").Value = "Male"

# --- Row 4: abed ak ---
$ws.Range("A4").Value = 123435234
$ws.Range("B4").Value = "abed"
$ws.Range("C4").Value = "ak"
$ws.Range("D4").Value = 24
$ws.Range("E4").Value = 70
$ws.Range("F4").Value = 176
$ws.Range("G4").Value = 524111123
$ws.Range("H4").Value = "A"
$ws.Range("I4").Value = "Male"

# --- Workbook-level defined name "ID" -> Sheet1!$A:$A ---
$wb.Names.Add('ID', 'Sheet1!$A:$A')

# --- Update selection to F1 (matches the saved cursor position) ---
[void]$ws.Range("F1").Select()
